$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension / used-range related header additions and shifted data

$ws.Cells.Item(1, 2).Value = 0
$ws.Cells.Item(1, 3).Value = 1
$ws.Cells.Item(1, 4).Value = 2
$ws.Cells.Item(1, 5).Value = 3
$ws.Cells.Item(1, 6).Value = 4
$ws.Cells.Item(1, 7).Value = 5
$ws.Cells.Item(1, 8).Value = 6
$ws.Cells.Item(1, 9).Value = 7
$ws.Cells.Item(1, 10).Value = 8
$ws.Cells.Item(1, 11).Value = 9
$ws.Cells.Item(1, 12).Value = 10
$ws.Cells.Item(1, 13).Value = 11
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = '100'
$ws.Cells.Item(2, 3).Value = '01'
$ws.Cells.Item(2, 4).Value = '18312'
$ws.Cells.Item(2, 5).Value = '4'
$ws.Cells.Item(2, 6).Value = 'F'
$ws.Cells.Item(2, 7).Value = 'MWF'
$ws.Cells.Item(2, 8).Value = 'Intro to Comparative Amer Stud'
$ws.Cells.Item(2, 9).Value = '0900'
$ws.Cells.Item(2, 10).Value = '0950am'
$ws.Cells.Item(2, 11).Value = 'KING'
$ws.Cells.Item(2, 12).Value = '123'
$ws.Cells.Item(2, 13).Value = 'Lee Shelley'
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = '101'
$ws.Cells.Item(3, 3).Value = '01'
$ws.Cells.Item(3, 4).Value = '18679'
$ws.Cells.Item(3, 5).Value = '4'
$ws.Cells.Item(3, 6).Value = 'F'
$ws.Cells.Item(3, 7).Value = 'TR'
$ws.Cells.Item(3, 8).Value = 'Intro to GSFS'
$ws.Cells.Item(3, 9).Value = '0130'
$ws.Cells.Item(3, 10).Value = '0250pm'
$ws.Cells.Item(3, 11).Value = 'KING'
$ws.Cells.Item(3, 12).Value = '123'
$ws.Cells.Item(3, 13).Value = 'LaGrotteria Angela'
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = '201'
$ws.Cells.Item(4, 3).Value = '01'
$ws.Cells.Item(4, 4).Value = '18313'
$ws.Cells.Item(4, 5).Value = '4'
$ws.Cells.Item(4, 6).Value = 'F'
$ws.Cells.Item(4, 7).Value = 'MWF'
$ws.Cells.Item(4, 8).Value = 'Latinas/os Comparative Perspec'
$ws.Cells.Item(4, 9).Value = '1000'
$ws.Cells.Item(4, 10).Value = '1050am'
$ws.Cells.Item(4, 11).Value = 'KING'
$ws.Cells.Item(4, 12).Value = '341'
$ws.Cells.Item(4, 13).Value = 'Perez Gina'
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = '264'
$ws.Cells.Item(5, 3).Value = '01'
$ws.Cells.Item(5, 4).Value = '19228'
$ws.Cells.Item(5, 5).Value = '4'
$ws.Cells.Item(5, 6).Value = 'F'
$ws.Cells.Item(5, 7).Value = 'TR'
$ws.Cells.Item(5, 8).Value = 'Abortion and Religion'
$ws.Cells.Item(5, 9).Value = '0130'
$ws.Cells.Item(5, 10).Value = '0245pm'
$ws.Cells.Item(5, 11).Value = 'KING'
$ws.Cells.Item(5, 12).Value = '327'
$ws.Cells.Item(5, 13).Value = 'Kamitsuka Margaret'
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = '278'
$ws.Cells.Item(6, 3).Value = '01'
$ws.Cells.Item(6, 4).Value = '19412'
$ws.Cells.Item(6, 5).Value = '4'
$ws.Cells.Item(6, 6).Value = 'F'
$ws.Cells.Item(6, 7).Value = 'TR'
$ws.Cells.Item(6, 8).Value = 'Gender, Race & War on Terror'
$ws.Cells.Item(6, 9).Value = '0300'
$ws.Cells.Item(6, 10).Value = '0415pm'
$ws.Cells.Item(6, 11).Value = 'KING'
$ws.Cells.Item(6, 12).Value = '241'
$ws.Cells.Item(6, 13).Value = 'Miller Kathryn'
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = '301'
$ws.Cells.Item(7, 3).Value = '01'
$ws.Cells.Item(7, 4).Value = '19422'
$ws.Cells.Item(7, 5).Value = '4'
$ws.Cells.Item(7, 6).Value = 'F'
$ws.Cells.Item(7, 7).Value = 'TR'
$ws.Cells.Item(7, 8).Value = 'Feminist Theory'
$ws.Cells.Item(7, 9).Value = '0835'
$ws.Cells.Item(7, 10).Value = '0950am'
$ws.Cells.Item(7, 11).Value = 'KING'
$ws.Cells.Item(7, 12).Value = '325'
$ws.Cells.Item(7, 13).Value = 'LaGrotteria Angela'
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = '304'
$ws.Cells.Item(8, 3).Value = '01'
$ws.Cells.Item(8, 4).Value = '19227'
$ws.Cells.Item(8, 5).Value = '4'
$ws.Cells.Item(8, 6).Value = 'F'
$ws.Cells.Item(8, 7).Value = 'TR'
$ws.Cells.Item(8, 8).Value = 'Transnational Feminisms'
$ws.Cells.Item(8, 9).Value = '1100'
$ws.Cells.Item(8, 10).Value = '1215pm'
$ws.Cells.Item(8, 11).Value = 'KING'
$ws.Cells.Item(8, 12).Value = '235'
$ws.Cells.Item(8, 13).Value = 'Miller Kathryn'
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = '306'
$ws.Cells.Item(9, 3).Value = '01'
$ws.Cells.Item(9, 4).Value = '18635'
$ws.Cells.Item(9, 5).Value = '4'
$ws.Cells.Item(9, 6).Value = 'F'
$ws.Cells.Item(9, 7).Value = 'W'
$ws.Cells.Item(9, 8).Value = 'Gender and Migration'
$ws.Cells.Item(9, 9).Value = '0230'
$ws.Cells.Item(9, 10).Value = '0420pm'
$ws.Cells.Item(9, 11).Value = 'RICE'
$ws.Cells.Item(9, 12).Value = '100B'
$ws.Cells.Item(9, 13).Value = 'Miller Kathryn'
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = '319'
$ws.Cells.Item(10, 3).Value = '01'
$ws.Cells.Item(10, 4).Value = '19337'
$ws.Cells.Item(10, 5).Value = '4'
$ws.Cells.Item(10, 6).Value = 'F'
$ws.Cells.Item(10, 7).Value = 'TR'
$ws.Cells.Item(10, 8).Value = 'Sexual ?Absences?'
$ws.Cells.Item(10, 9).Value = '0300'
$ws.Cells.Item(10, 10).Value = '0415pm'
$ws.Cells.Item(10, 11).Value = 'KING'
$ws.Cells.Item(10, 12).Value = '325'
$ws.Cells.Item(10, 13).Value = 'Cerankowski KJ'
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = '330'
$ws.Cells.Item(11, 3).Value = '01'
$ws.Cells.Item(11, 4).Value = '19194'
$ws.Cells.Item(11, 5).Value = '4'
$ws.Cells.Item(11, 6).Value = 'F'
$ws.Cells.Item(11, 7).Value = 'M'
$ws.Cells.Item(11, 8).Value = 'Relg Gender Sexuality in India'
$ws.Cells.Item(11, 9).Value = '0230'
$ws.Cells.Item(11, 10).Value = '0420pm'
$ws.Cells.Item(11, 11).Value = 'AJLC'
$ws.Cells.Item(11, 12).Value = '102A'
$ws.Cells.Item(11, 13).Value = 'Bachrach Emilia'
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = '347'
$ws.Cells.Item(12, 3).Value = '01'
$ws.Cells.Item(12, 4).Value = '19420'
$ws.Cells.Item(12, 5).Value = '4'
$ws.Cells.Item(12, 6).Value = 'F'
$ws.Cells.Item(12, 7).Value = 'TR'
$ws.Cells.Item(12, 8).Value = 'Queer Postwar New York & Paris'
$ws.Cells.Item(12, 9).Value = '0930'
$ws.Cells.Item(12, 10).Value = '1050am'
$ws.Cells.Item(12, 11).Value = 'PETE'
$ws.Cells.Item(12, 12).Value = '232'
$ws.Cells.Item(12, 13).Value = 'O''Connor Patrick'
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = '400'
$ws.Cells.Item(13, 3).Value = '01'
$ws.Cells.Item(13, 4).Value = '14677'
$ws.Cells.Item(13, 5).Value = '0'
$ws.Cells.Item(13, 6).Value = 'S'
$ws.Cells.Item(13, 8).Value = 'Senior'
$ws.Cells.Item(13, 9).Value = 'Capstone'
$ws.Cells.Item(13, 10).Value = 'TBA'
$ws.Cells.Item(13, 11).Value = 'TBA'
$ws.Cells.Item(13, 12).Value = 'Mattson Greggor'
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = '500'
$ws.Cells.Item(14, 3).Value = '01'
$ws.Cells.Item(14, 4).Value = '14560'
$ws.Cells.Item(14, 5).Value = '4'
$ws.Cells.Item(14, 7).Value = 'F'
$ws.Cells.Item(14, 8).Value = 'Honors'
$ws.Cells.Item(14, 9).Value = 'TBA'
$ws.Cells.Item(14, 10).Value = 'TBA'
$ws.Cells.Item(14, 11).Value = 'Mattson Greggor'
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = '995F'
$ws.Cells.Item(15, 3).Value = '01'
$ws.Cells.Item(15, 4).Value = '18073'
$ws.Cells.Item(15, 5).Value = '4'
$ws.Cells.Item(15, 6).Value = 'F'
$ws.Cells.Item(15, 8).Value = 'Private Reading'
$ws.Cells.Item(15, 10).Value = 'Full'
$ws.Cells.Item(15, 11).Value = 'TBA'
$ws.Cells.Item(15, 12).Value = 'TBA'
$ws.Cells.Item(15, 13).Value = 'Mattson Greggor'
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = '995H'
$ws.Cells.Item(16, 3).Value = '01'
$ws.Cells.Item(16, 4).Value = '18074'
$ws.Cells.Item(16, 5).Value = '2'
$ws.Cells.Item(16, 6).Value = 'F'
$ws.Cells.Item(16, 8).Value = 'Private Reading'
$ws.Cells.Item(16, 10).Value = 'Half'
$ws.Cells.Item(16, 11).Value = 'TBA'
$ws.Cells.Item(16, 12).Value = 'TBA'
$ws.Cells.Item(16, 13).Value = 'Mattson Greggor'

# Clear cells that no longer hold data after the column reshuffle
$ws.Cells.Item(13, 7).ClearContents()
$ws.Cells.Item(14, 6).ClearContents()
$ws.Cells.Item(15, 7).ClearContents()
$ws.Cells.Item(15, 9).ClearContents()
$ws.Cells.Item(16, 7).ClearContents()
$ws.Cells.Item(16, 9).ClearContents()
